# Populate the new "LLBV3 Header / Function" column (E) with a short note
# on what each populated pin is used for / connected to, matching the
# current state of the board. Entered in the same order the author typed
# them so newly-introduced strings land in the shared string table the
# same way a live Excel edit session would produce them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "LLBV3 Header" -> "LLBV3 Header / Function"
$ws.Range("E1").Value = "LLBV3 Header / Function"

$ws.Range("E2").Value  = "MCP 2515 interrupt on received frames"
$ws.Range("E3").Value  = "USB Serial"
$ws.Range("E4").Value  = "USB Serial"
$ws.Range("E6").Value  = "X3, for power on board"
$ws.Range("E7").Value  = "X3, for power on board"
$ws.Range("E8").Value  = "E-stop jumper, also X3"
$ws.Range("E17").Value = "Steering header"
$ws.Range("E18").Value = "X3, for power on board"
$ws.Range("E19").Value = "wheel hall switch header"

# These two share similar text, entered with the "all SPI devices" note
# first (rows 21-23), then the SPI_SLAVE-specific note on row 20.
$ws.Range("E21").Value = "all SPI devices, SPI header"
$ws.Range("E22").Value = "all SPI devices, SPI header"
$ws.Range("E23").Value = "all SPI devices, SPI header"
$ws.Range("E20").Value = "SPI header (this pin tells the mega to be a slave)"

$ws.Range("E24").Value = "X3, for power on board"
$ws.Range("E25").Value = "X3, for power on board"
$ws.Range("E26").Value = "X3, for power on board"
$ws.Range("E27").Value = "X3, for power on board"
$ws.Range("E36").Value = "MCP2515 slave selection"
$ws.Range("E37").Value = "DAC slave selection"
$ws.Range("E52").Value = "on-board relay"
$ws.Range("E53").Value = "on-board buzzer"
$ws.Range("E54").Value = "on-board relay"
$ws.Range("E55").Value = "X3, no purpose assgined"
$ws.Range("E57").Value = "X3, for power on board"
$ws.Range("E64").Value = "X3, for power on board"
$ws.Range("E65").Value = "X3, for power on board"
$ws.Range("E77").Value = "X3, for power on board"
$ws.Range("E79").Value = "X3, for power on board"
$ws.Range("E91").Value = "Steering header"
$ws.Range("E92").Value = "Steering header"
$ws.Range("E95").Value = "Steering header"
$ws.Range("E96").Value = "Steering header"

# Selection state as left by the author after the edit.
$ws.Range("E58").Select()
